$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 17 data (Trasporte / transport / Transport)
$ws.Range("A17").Value = "Trasporte"
$ws.Range("B17").Value = "transport"
$ws.Range("C17").Value = "Transport"
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = "circle-o"
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0

# Update the selected cell, to match the diff (active cell D17)
$ws.Range("D17").Select()
